# Updating BOM for mirrored STL
# Insert a new BOM line for the mirrored "FRAME VERTEX WITH FOOT" printed
# part (part # 1010) directly below the existing "FRAME VERTEX WITH FOOT"
# line (row 86), and halve the quantity of the original (non-mirrored)
# part from 4 down to 2 since two of the four are now printed mirrored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 87:98 down to 88:99, leaving a blank row 87 for the new part.
$ws.Rows.Item(87).Insert()

# Existing "FRAME VERTEX WITH FOOT" (row 86) quantity drops from 4 to 2.
$ws.Cells.Item(86, 5).Value = 2

# New row 87: mirrored variant of the frame vertex foot.
$ws.Cells.Item(87, 1).Value = 1010
$ws.Cells.Item(87, 2).Value = "00-000"
$ws.Cells.Item(87, 3).Value = "FRAME VERTEX WITH FOOT(MIRRORED)"
$ws.Cells.Item(87, 4).Value = "PRINTED PARTS"
$ws.Cells.Item(87, 5).Value = 2

# Both the original and the new mirrored-part row use the tighter,
# auto-fit row height seen once the two rows are grouped together.
$ws.Rows.Item(86).RowHeight = 13.8
$ws.Rows.Item(87).RowHeight = 13.8

# Leave the selection near the newly inserted row, matching where the
# edit was made.
[void]$ws.Range("A55").Select()
[void]$ws.Range("C87").Select()
